$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before "Login" (column F), shifting Login to column G
$ws.Columns("F").Insert()

# Set header for new column
$ws.Range("F1").Value = "Practical work"

# Set the practical-work values for each student row
$ws.Range("F2").Value = "PW1"
$ws.Range("F3").Value = "PW2"
$ws.Range("F4").Value = "PW3"
$ws.Range("F5").Value = "PW4"
$ws.Range("F6").Value = "PW1"
$ws.Range("F7").Value = "PW2"
$ws.Range("F8").Value = "PW3"
$ws.Range("F9").Value = "PW4"

$ws.Range("F10").Select()
